# Auto-generated Excel COM-interop edit script
# Applies the symbol-list update described in the commit:
# "Updated symbol list on Sun Jan  1 10:57:07 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (preserves exact formatting such as
# trailing zeros "3.020" / leading zeros / percent signs) instead of letting
# Excel auto-convert numeric-looking or percent-looking strings.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Row 2
Set-TextValue "D2" "243.98"
Set-TextValue "E2" "-0.81%"

# Row 3
Set-TextValue "D3" "26.99"
Set-TextValue "E3" "3.38%"

# Row 4
Set-TextValue "D4" "5.156"
Set-TextValue "E4" "1.04%"

# Row 5
Set-TextValue "D5" "0.05628"
Set-TextValue "E5" "0.42%"

# Row 6
Set-TextValue "D6" "6.477"

# Row 7
Set-TextValue "D7" "0.8166"
Set-TextValue "E7" "0.67%"

# Row 8
Set-TextValue "D8" "0.8315"
Set-TextValue "E8" "-1.70%"

# Row 9
Set-TextValue "D9" "0.1327"
Set-TextValue "E9" "-1.05%"

# Row 10
Set-TextValue "D10" "0.06926"
Set-TextValue "E10" "-0.39%"

# Row 11
Set-TextValue "D11" "0.02916"
Set-TextValue "E11" "2.39%"

# Row 12
Set-TextValue "E12" "0.00%"

# Row 13
Set-TextValue "D13" "0.001518"
Set-TextValue "E13" "0.29%"

# Row 14
Set-TextValue "D14" "0.04225"
Set-TextValue "E14" "-10.21%"

# Row 15
Set-TextValue "D15" "0.0005988"
Set-TextValue "E15" "-0.34%"

# Row 16
Set-TextValue "E16" "-0.02%"

# Row 17
Set-TextValue "D17" "3.594"
Set-TextValue "E17" "0.89%"

# Row 18
Set-TextValue "D18" "3.020"
Set-TextValue "E18" "-0.04%"

# Row 19
Set-TextValue "D19" "2.226"
Set-TextValue "E19" "5.08%"

# Row 21
Set-TextValue "D21" "0.03092"
Set-TextValue "E21" "-3.68%"

# Row 22
Set-TextValue "E22" "-2.13%"

# Row 23
Set-TextValue "D23" "3.745"
Set-TextValue "E23" "-0.15%"

# Row 24
Set-TextValue "D24" "0.1374"

# Row 25
Set-TextValue "E25" "-1.65%"

# Row 26
Set-TextValue "D26" "0.004484"
Set-TextValue "E26" "-2.86%"

# Row 27
Set-TextValue "D27" "0.00009795"
Set-TextValue "E27" "2.05%"

# Row 28
Set-TextValue "E28" "-0.45%"

# Row 40
Set-TextValue "D40" "0.03651"
Set-TextValue "E40" "-0.18%"

# Row 41
Set-TextValue "B41" "BKEXToken"
Set-TextValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1376"
Set-TextValue "E41" "1.81%"

# Row 42
Set-TextValue "B42" "CEJI"
Set-TextValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002619"
Set-TextValue "E42" "-1.53%"

# Row 43
Set-TextValue "B43" "KickToken"
Set-TextValue "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003432"
Set-TextValue "E43" "-44.00%"

# Row 44
Set-TextValue "D44" "0.008172"
Set-TextValue "E44" "5.32%"

# Row 45
Set-TextValue "D45" "0.00005341"
Set-TextValue "E45" "0.93%"

# Row 46
Set-TextValue "E46" "0.00%"

# Row 47
Set-TextValue "E47" "-18.05%"

# Row 48
Set-TextValue "D48" "0.002640"
Set-TextValue "E48" "28.84%"

# Row 49
Set-TextValue "E49" "0.00%"

# Row 50
Set-TextValue "E50" "0.00%"
